$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Journal sheet (sheet1)
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Journal")
$lo1 = $ws1.ListObjects.Item("Tableau1")

# New "Commentaire" header/column for the table (set before the other new
# strings so the shared-string table keeps the same allocation order as
# the authored workbook: Commentaire, comment text, then the two export
# descriptions).
$ws1.Cells.Item(1, 6).Value2 = "Commentaire"

# Update the 3 duration cells that were corrected (0:45->1:00, 1:00->1:10,
# 1:20->1:25).
$ws1.Cells.Item(28, 3).Value2 = 0.041666666666666664
$ws1.Cells.Item(30, 3).Value2 = 0.048611111111111112
$ws1.Cells.Item(31, 3).Value2 = 0.059027777777777783

# Row 32: clone formatting from row 31, then overwrite with the new entry.
$ws1.Range("A31:F31").Copy() | Out-Null
$ws1.Range("A32:F32").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(32, 1).Value2 = 44985
$ws1.Cells.Item(32, 2).Value2 = 4
$ws1.Cells.Item(32, 3).Value2 = 0.03125
$ws1.Cells.Item(32, 4).Value2 = "Documentation"
$ws1.Cells.Item(32, 5).Value2 = "Rédaction des futures questions à poser pour les entretiens"
$ws1.Cells.Item(32, 6).Value2 = "J'ai rencontré beaucoup de difficulté à formuler et a trouver des questions pour ce document"

# Row 33: clone formatting from row 14 (an existing "Export" row), 5 columns only.
$ws1.Range("A14:E14").Copy() | Out-Null
$ws1.Range("A33:E33").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(33, 1).Value2 = 44985
$ws1.Cells.Item(33, 2).Value2 = 4
$ws1.Cells.Item(33, 3).Value2 = 0.017361111111111112
$ws1.Cells.Item(33, 4).Value2 = "Export"
$ws1.Cells.Item(33, 5).Value2 = "Export des différents types de composant dans un document Word"

# Row 34: same shape as row 33.
$ws1.Range("A14:E14").Copy() | Out-Null
$ws1.Range("A34:E34").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(34, 1).Value2 = 44985
$ws1.Cells.Item(34, 2).Value2 = 4
$ws1.Cells.Item(34, 3).Value2 = 0.013888888888888888
$ws1.Cells.Item(34, 4).Value2 = "Export"
$ws1.Cells.Item(34, 5).Value2 = "Export des différents logiciels dans un document Word"

# Grow the table to cover the new column + new rows, then fix up the new
# column's header name (ListObjects resolve the column name from the
# header cell, so this must happen after the resize).
$lo1.Resize($ws1.Range("A1:F34"))
$ws1.Cells.Item(1, 6).Value2 = "Commentaire"

# Give the new header cell the same look as the other table headers
# (centered, slightly larger font).
$ws1.Cells.Item(1, 6).Font.Size = 12
$ws1.Cells.Item(1, 6).HorizontalAlignment = -4108

$ws1.Range("C31").Select()

# ----------------------------------------------------------------------
# Totaux sheet (sheet2)
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Totaux")
$lo2 = $ws2.ListObjects.Item("Tableau2")

# Old totals row (12) becomes row 13; clone its formatting down first.
$ws2.Range("A12:B12").Copy() | Out-Null
$ws2.Range("A13:B13").PasteSpecial(-4122) | Out-Null

# Row 11 gains its weekly-total formula (it already has the date).
$ws2.Cells.Item(11, 2).Formula = "=SUM(Journal!C27:C31)"
$ws2.Cells.Item(11, 2).NumberFormat = $ws2.Cells.Item(10, 2).NumberFormat

# Row 12 becomes a normal data row for the new week.
$ws2.Range("A10:B10").Copy() | Out-Null
$ws2.Range("A12:B12").PasteSpecial(-4122) | Out-Null
$ws2.Cells.Item(12, 1).Value2 = 44985
$ws2.Cells.Item(12, 2).Formula = "=SUM(Journal!C32:C34)"

# Row 13: restore the grand-total label + formula over the new range.
$ws2.Cells.Item(13, 1).Value2 = "Total"
$ws2.Cells.Item(13, 2).Formula = "=SUM(B2:B12)"

$lo2.Resize($ws2.Range("A1:B13"))

$ws2.Range("F15").Select()
$ws2.Activate()

$excel.Calculate()
